$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 8: eu_sub subsidy entry
$ws.Range("A8").Value = "eu_sub"
$ws.Range("B8").Value = 150
$ws.Range("C8").Value = "NA"
$ws.Range("D8").Value = 150
$ws.Range("E8").Value = "const"
$ws.Range("F8").Value = "€/ha"

# Update the selected cell to match the saved selection state
$ws.Range("D10").Select()
